# Scheduled-runner market-data refresh for the Leviathan Profits workbook.
# Updates the computed market-price / profit columns (H-N) on each job sheet
# with freshly pulled currentAveragePrice figures; leve/recipe metadata (A-G)
# is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 827.7143
$ws.Range("I12").Value = 658.4
$ws.Range("K12").Value = 658.4
$ws.Range("M12").Value = -488.4
$ws.Range("H15").Value = 1290.7142
$ws.Range("I15").Value = 1290.7142
$ws.Range("K15").Value = 3872.1426
$ws.Range("M15").Value = -3703.1426
$ws.Range("H74").Value = 4470.4287
$ws.Range("I74").Value = 3439.3333
$ws.Range("K74").Value = 3439.3333
$ws.Range("M74").Value = -2503.3333
$ws.Range("H77").Value = 4470.4287
$ws.Range("I77").Value = 3439.3333
$ws.Range("K77").Value = 17196.6665
$ws.Range("M77").Value = -12516.6665
$ws.Range("H87").Value = 33354
$ws.Range("J87").Value = 33354
$ws.Range("L87").Value = 33354
$ws.Range("N87").Value = -35850
$ws.Range("H90").Value = 33354
$ws.Range("J90").Value = 33354
$ws.Range("L90").Value = 100062
$ws.Range("N90").Value = -112542
$ws.Range("H132").Value = 3614.0952
$ws.Range("I132").Value = 1694.037
$ws.Range("J132").Value = 7070.2
$ws.Range("K132").Value = 5082.111
$ws.Range("L132").Value = 21210.6
$ws.Range("M132").Value = -2552.111
$ws.Range("N132").Value = -26270.6
$ws.Range("H137").Value = 29261.555
$ws.Range("I137").Value = 1575.8
$ws.Range("J137").Value = 92183.73
$ws.Range("K137").Value = 4727.4
$ws.Range("L137").Value = 276551.19
$ws.Range("M137").Value = -2177.4
$ws.Range("N137").Value = -281651.19
$ws.Range("H138").Value = 2040.0577
$ws.Range("I138").Value = 1365.9166
$ws.Range("J138").Value = 2617.8928
$ws.Range("K138").Value = 4097.7498
$ws.Range("L138").Value = 7853.678400000001
$ws.Range("M138").Value = 1042.2502
$ws.Range("N138").Value = -18133.6784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2356.5715
$ws.Range("I63").Value = 2353.3076
$ws.Range("K63").Value = 2353.3076
$ws.Range("M63").Value = -1667.3076
$ws.Range("H66").Value = 2356.5715
$ws.Range("I66").Value = 2353.3076
$ws.Range("K66").Value = 11766.538
$ws.Range("M66").Value = -8334.538
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 33883.79
$ws.Range("I132").Value = 33883.79
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 101651.37
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -99121.37
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5417.9165
$ws.Range("J20").Value = 7748.222
$ws.Range("L20").Value = 7748.222
$ws.Range("N20").Value = -8242.222
$ws.Range("H86").Value = 2034.3684
$ws.Range("I86").Value = 1822.5454
$ws.Range("J86").Value = 2325.625
$ws.Range("K86").Value = 1822.5454
$ws.Range("L86").Value = 2325.625
$ws.Range("M86").Value = -699.5454
$ws.Range("N86").Value = -4571.625
$ws.Range("H89").Value = 2034.3684
$ws.Range("I89").Value = 1822.5454
$ws.Range("J89").Value = 2325.625
$ws.Range("K89").Value = 9112.726999999999
$ws.Range("L89").Value = 11628.125
$ws.Range("M89").Value = -3496.726999999999
$ws.Range("N89").Value = -22860.125
$ws.Range("H105").Value = 3229741.8
$ws.Range("I105").Value = 3850306.2
$ws.Range("K105").Value = 3850306.2
$ws.Range("M105").Value = -3848559.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1585.931
$ws.Range("I31").Value = 1615.76
$ws.Range("K31").Value = 1615.76
$ws.Range("M31").Value = -1320.76
$ws.Range("H34").Value = 1585.931
$ws.Range("I34").Value = 1615.76
$ws.Range("K34").Value = 1615.76
$ws.Range("M34").Value = -1413.76
$ws.Range("H58").Value = 2119.2222
$ws.Range("I58").Value = 2179.8333
$ws.Range("J58").Value = 1998
$ws.Range("K58").Value = 2179.8333
$ws.Range("L58").Value = 1998
$ws.Range("M58").Value = -1976.8333
$ws.Range("N58").Value = -2404
$ws.Range("H86").Value = 9666.333000000001
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 9666.333000000001
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 9666.333000000001
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -11912.333
$ws.Range("H89").Value = 9666.333000000001
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 9666.333000000001
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 48331.665
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -59563.665
$ws.Range("H132").Value = 3686
$ws.Range("J132").Value = 3948
$ws.Range("L132").Value = 11844
$ws.Range("N132").Value = -16904
$ws.Range("H136").Value = 2119.2222
$ws.Range("I136").Value = 2179.8333
$ws.Range("J136").Value = 1998
$ws.Range("K136").Value = 6539.499899999999
$ws.Range("L136").Value = 5994
$ws.Range("M136").Value = -3989.499899999999
$ws.Range("N136").Value = -11094

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 9662.210999999999
$ws.Range("I56").Value = 9662.210999999999
$ws.Range("K56").Value = 9662.210999999999
$ws.Range("M56").Value = -9132.210999999999
$ws.Range("H113").Value = 573.5714
$ws.Range("I113").Value = 498.66666
$ws.Range("J113").Value = 629.75
$ws.Range("K113").Value = 1495.99998
$ws.Range("L113").Value = 1889.25
$ws.Range("M113").Value = 674.0000199999999
$ws.Range("N113").Value = -6229.25
$ws.Range("H141").Value = 3125.25
$ws.Range("I141").Value = 3074
$ws.Range("K141").Value = 9222
$ws.Range("M141").Value = -4042

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 31023.273
$ws.Range("I97").Value = 39348
$ws.Range("K97").Value = 39348
$ws.Range("M97").Value = -38852

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1240
$ws.Range("I22").Value = 1066.6666
$ws.Range("K22").Value = 1066.6666
$ws.Range("M22").Value = -771.6666
$ws.Range("H27").Value = 1240
$ws.Range("I27").Value = 1066.6666
$ws.Range("K27").Value = 1066.6666
$ws.Range("M27").Value = -959.6666
$ws.Range("H68").Value = 2191.2222
$ws.Range("I68").Value = 1866.8334
$ws.Range("J68").Value = 2840
$ws.Range("K68").Value = 1866.8334
$ws.Range("L68").Value = 2840
$ws.Range("M68").Value = -1117.8334
$ws.Range("N68").Value = -4338
$ws.Range("H71").Value = 2191.2222
$ws.Range("I71").Value = 1866.8334
$ws.Range("J71").Value = 2840
$ws.Range("K71").Value = 9334.166999999999
$ws.Range("L71").Value = 14200
$ws.Range("M71").Value = -5590.166999999999
$ws.Range("N71").Value = -21688
$ws.Range("H82").Value = 1274.5
$ws.Range("I82").Value = 900
$ws.Range("J82").Value = 1399.3334
$ws.Range("K82").Value = 900
$ws.Range("L82").Value = 1399.3334
$ws.Range("M82").Value = -539
$ws.Range("N82").Value = -2121.3334
$ws.Range("H85").Value = 1274.5
$ws.Range("I85").Value = 900
$ws.Range("J85").Value = 1399.3334
$ws.Range("K85").Value = 900
$ws.Range("L85").Value = 1399.3334
$ws.Range("M85").Value = 348
$ws.Range("N85").Value = -3895.3334
$ws.Range("H87").Value = 20189
$ws.Range("J87").Value = 20189
$ws.Range("L87").Value = 20189
$ws.Range("N87").Value = -22435
$ws.Range("H90").Value = 20189
$ws.Range("J90").Value = 20189
$ws.Range("L90").Value = 60567
$ws.Range("N90").Value = -71799
$ws.Range("H93").Value = 2192
$ws.Range("I93").Value = 1888.4546
$ws.Range("K93").Value = 1888.4546
$ws.Range("M93").Value = -640.4546
$ws.Range("H122").Value = 10100.134
$ws.Range("I122").Value = 11375.25
$ws.Range("J122").Value = 4999.6665
$ws.Range("K122").Value = 34125.75
$ws.Range("L122").Value = 14998.9995
$ws.Range("M122").Value = -31675.75
$ws.Range("N122").Value = -19898.9995
$ws.Range("H132").Value = 2541.7368
$ws.Range("I132").Value = 2294.111
$ws.Range("J132").Value = 6999
$ws.Range("K132").Value = 6882.333
$ws.Range("L132").Value = 20997
$ws.Range("M132").Value = -4352.333
$ws.Range("N132").Value = -26057
$ws.Range("H136").Value = 2824
$ws.Range("I136").Value = 2086.4167
$ws.Range("K136").Value = 6259.250100000001
$ws.Range("M136").Value = -3709.250100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9652.166999999999
$ws.Range("I2").Value = 12503.25
$ws.Range("J2").Value = 3950
$ws.Range("K2").Value = 12503.25
$ws.Range("L2").Value = 3950
$ws.Range("M2").Value = -12391.25
$ws.Range("N2").Value = -4174
$ws.Range("H62").Value = 7618.8335
$ws.Range("I62").Value = 8011.923
$ws.Range("J62").Value = 6596.8
$ws.Range("K62").Value = 8011.923
$ws.Range("L62").Value = 6596.8
$ws.Range("M62").Value = -7387.923
$ws.Range("N62").Value = -7844.8
$ws.Range("H65").Value = 7618.8335
$ws.Range("I65").Value = 8011.923
$ws.Range("J65").Value = 6596.8
$ws.Range("K65").Value = 40059.615
$ws.Range("L65").Value = 32984
$ws.Range("M65").Value = -36939.615
$ws.Range("N65").Value = -39224
$ws.Range("H100").Value = 1986.45
$ws.Range("I100").Value = 2437.6
$ws.Range("K100").Value = 4875.2
$ws.Range("M100").Value = -4334.2
$ws.Range("H105").Value = 40615
$ws.Range("J105").Value = 40615
$ws.Range("L105").Value = 40615
$ws.Range("N105").Value = -47603
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800
$ws.Range("H136").Value = 1050.3684
$ws.Range("I136").Value = 997
$ws.Range("K136").Value = 2991
$ws.Range("M136").Value = -441
